$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.899.78'
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.287.66'
$ws.Range("E3").Value = '  -3.47%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.32'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.56'
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -3.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.25'
$ws.Range("E10").Value = '  -6.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("E12").Value = '  -4.10%  '
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("E14").Value = '  -5.16%  '
$ws.Range("E15").Value = '  -6.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.634.39'
$ws.Range("E16").Value = '  -3.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.288.08'
$ws.Range("E17").Value = '  -3.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.867.56'
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("E19").Value = '  -2.90%  '
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '282.78'
$ws.Range("E21").Value = '  +9.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.61'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.15'
$ws.Range("E23").Value = '  -4.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.00'
$ws.Range("E24").Value = '  +6.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  -4.18%  '
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.75'
$ws.Range("E27").Value = '  -6.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.39'
$ws.Range("E28").Value = '  +6.20%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.06'
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.02'
$ws.Range("E31").Value = '  -4.75%  '
$ws.Range("E32").Value = '  -3.09%  '
$ws.Range("E33").Value = '  -3.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.81'
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.135'
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  -7.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.53'
$ws.Range("E37").Value = '  -4.69%  '
$ws.Range("E38").Value = '  +8.38%  '
$ws.Range("E39").Value = '  -4.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.68'
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.90'
$ws.Range("E41").Value = '  +13.27%  '
$ws.Range("E42").Value = '  -5.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.51'
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("E45").Value = '  -7.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '114.20'
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.91'
$ws.Range("E47").Value = '  -3.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.60'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("E50").Value = '  -5.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.550.29'
$ws.Range("E51").Value = '  -1.38%  '
